function Replace-Text($doc, $old, $new) {
    $doc.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

$d = $word.ActiveDocument

# Header date
Replace-Text $d "2024-04-10 Wednesday" "2024-04-11 Thursday"

# Simple one-for-one equation replacements (unique text across the document)
Replace-Text $d "36÷2=" "35÷6="
Replace-Text $d "29÷3=" "28÷6="
Replace-Text $d "61÷8=" "63÷4="
Replace-Text $d "90÷7=" "40÷2="
Replace-Text $d "71÷5=" "26÷5="
Replace-Text $d "39÷5=" "69÷3="
Replace-Text $d "31÷5=" "20÷6="
Replace-Text $d "66÷5=" "46÷3="
Replace-Text $d "44÷7=" "19÷5="
Replace-Text $d "10÷6=" "12÷8="
Replace-Text $d "48÷4=" "34÷2="
Replace-Text $d "95÷2=" "76÷5="
Replace-Text $d "27÷9=" "50÷7="
Replace-Text $d "68÷8=" "55÷4="
Replace-Text $d "10÷8=" "17÷6="
Replace-Text $d "60÷8=" "93÷2="
Replace-Text $d "51÷8=" "99÷5="
Replace-Text $d "45÷5=" "60÷3="
Replace-Text $d "69÷7=" "49÷4="
Replace-Text $d "16÷3=" "13÷9="

# Last table row: 54÷9=, 79÷9=, 36÷5=, 40÷3=, 66÷7= (5 cells)
# becomes: 82÷3=, 36÷5=, 53÷4=, 67÷3=, 32÷3= (still 5 cells) -- the line-based
# xml diff made this look like a cell removal + cell insertion, but the
# table keeps the same cell count; only the text content shifts.
$table = $d.Tables.Item(1)
$row = $table.Rows.Item(17)
$row.Cells.Item(1).Range.Text = "82÷3="
$row.Cells.Item(2).Range.Text = "36÷5="
$row.Cells.Item(3).Range.Text = "53÷4="
$row.Cells.Item(4).Range.Text = "67÷3="
$row.Cells.Item(5).Range.Text = "32÷3="
